# Update cryptocurrency price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we are about to rewrite are treated as
# plain text (matching the workbook's inline-string storage) so that
# Excel does not silently reinterpret values like "215.82" as numbers.
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D32","D33","D34","D35","D37","D38","D41","D43","D45","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "81.708.59"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.148.00"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "215.82"
$ws.Range("E5").Value = "  +4.98%  "
$ws.Range("D6").Value = "615.63"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("D7").Value = "0.284"
$ws.Range("E7").Value = "  +14.74%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("D10").Value = "3.147.49"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").Value = "0.593"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("D12").Value = "0.0000253"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "5.27"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "3.739.43"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "31.83"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "81.753.39"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "3.151.78"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "3.19"
$ws.Range("E19").Value = "  +7.63%  "
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -5.39%  "
$ws.Range("D21").Value = "431.48"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "8.85"
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("D23").Value = "5.09"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").Value = "7.29"
$ws.Range("E24").Value = "  +6.85%  "
$ws.Range("D25").Value = "5.20"
$ws.Range("E25").Value = "  +7.76%  "
$ws.Range("D26").Value = "11.72"
$ws.Range("E26").Value = "  +7.49%  "
$ws.Range("D27").Value = "3.313.60"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "76.32"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "8.94"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("D33").Value = "564.11"
$ws.Range("E33").Value = "  +5.25%  "
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  +18.96%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").Value = "22.47"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +8.71%  "
$ws.Range("D41").Value = "0.402"
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("D43").Value = "2.99"
$ws.Range("E43").Value = "  +13.08%  "
$ws.Range("E44").Value = "  +8.53%  "
$ws.Range("D45").Value = "159.67"
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D47").Value = "185.58"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("D48").Value = "44.35"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").Value = "26.27"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("E51").Value = "  -5.88%  "

# Restore default (Normal) cell style on the price cells so only the
# underlying value changes are reflected, keeping formatting identical
# to the original workbook.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
